# Update the division-problem worksheet cells to the newly generated values.
# Each Find/Replace targets one specific cell's text; matching is case-
# sensitive and whole strings (including the trailing "=") so that only the
# intended run is affected and duplicates are avoided.

$d = $word.ActiveDocument

$replacements = @(
    @("81÷4=", "94÷3="),
    @("87÷4=", "47÷8="),
    @("92÷2=", "21÷7="),
    @("87÷7=", "34÷8="),
    @("64÷8=", "37÷8="),
    @("58÷7=", "80÷9="),
    @("51÷7=", "23÷4="),
    @("20÷9=", "44÷9="),
    @("49÷4=", "85÷5="),
    @("42÷5=", "35÷6="),
    @("75÷7=", "80÷6="),
    @("25÷2=", "97÷3="),
    @("51÷8=", "30÷3="),
    @("82÷2=", "28÷2="),
    @("48÷4=", "98÷9="),
    @("76÷9=", "51÷6="),
    @("98÷6=", "48÷3="),
    @("91÷7=", "76÷8="),
    @("14÷3=", "53÷5="),
    @("56÷7=", "84÷8="),
    @("77÷6=", "50÷4="),
    @("16÷5=", "24÷3="),
    @("30÷9=", "40÷4="),
    @("18÷8=", "53÷7=")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $range = $d.Content
    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}
